$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column E (copy formatting from the D1 header cell, then set text)
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "ERROR"

# Uncertainty values for each data row (row 2..28), aligned with column D spectroscopic factors
$errors = @(
    0.001719006577617426,
    0.05644954807469117,
    0.006012585757649582,
    0.0003458087934066839,
    0.009270130991122661,
    0.01600016777948418,
    0.01339535639708222,
    0.02022111690960823,
    0.004312818147205043,
    0.0006779208308500735,
    0.01578709551130119,
    0.0002906876987435252,
    0.002605447083208151,
    0.0002447672255863539,
    0.00560724076035283,
    0.002919866890431535,
    0.01265478670347784,
    0.002419325947724441,
    0.001271768185997996,
    0.002082889476641578,
    0.0009971253043785204,
    0.001927108398537859,
    0.0005541067359372513,
    0.00546960176799836,
    0.001101711701838224,
    0.007584759979328066,
    0.005182531202852024
)

for ($i = 0; $i -lt $errors.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 5).Value = $errors[$i]
}

# Row 21 (index 19) corrections: L column (C) and SPECTROSCOPIC_FACTOR (D)
$ws.Range("C21").Value = 3
$ws.Range("D21").Value = 0.02110661336330132
